$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.109.96"
$ws.Range("E2").Value = "  +3.55%  "
$ws.Range("D3").Value = "4.049.12"
$ws.Range("E3").Value = "  +3.17%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'519.69"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").Value = "'148.44"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "'0.738"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").Value = "'0.0000336"
$ws.Range("E11").Value = "  +0.53%  "
$ws.Range("D12").Value = "'47.05"
$ws.Range("E12").Value = "  +10.63%  "
$ws.Range("D13").Value = "'10.80"
$ws.Range("E13").Value = "  +4.77%  "
$ws.Range("D14").Value = "4.698.40"
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("D15").Value = "4.090.21"
$ws.Range("E15").Value = "  +3.96%  "
$ws.Range("D16").Value = "'21.33"
$ws.Range("E16").Value = "  +6.82%  "
$ws.Range("D17").Value = "'14.28"
$ws.Range("E17").Value = "  +2.52%  "
$ws.Range("E18").Value = "  -2.00%  "
$ws.Range("E19").Value = "  -1.73%  "
$ws.Range("D20").Value = "72.154.12"
$ws.Range("E20").Value = "  +3.77%  "
$ws.Range("D21").Value = "'443.73"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").Value = "'95.69"
$ws.Range("E22").Value = "  +8.60%  "
$ws.Range("D23").Value = "'3.52"
$ws.Range("E23").Value = "  +4.91%  "
$ws.Range("D24").Value = "'14.42"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "'12.15"
$ws.Range("E25").Value = "  +5.01%  "
$ws.Range("D26").Value = "'4.04"
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").Value = "'11.23"
$ws.Range("E27").Value = "  +4.71%  "
$ws.Range("D28").Value = "'37.06"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'3.11"
$ws.Range("E29").Value = "  +9.71%  "
$ws.Range("D30").Value = "'710.33"
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Value = "'13.48"
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.130"
$ws.Range("E32").Value = "  +2.41%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'6.92"
$ws.Range("E33").Value = "  +16.11%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "'67.66"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0911"
$ws.Range("E35").Value = "  +8.64%  "
$ws.Range("B36").Value = "TheGraph"
$ws.Range("C36").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D36").Value = "'0.442"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("B37").Value = "ThetaToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D37").Value = "'3.69"
$ws.Range("E37").Value = "  +24.62%  "
$ws.Range("D38").Value = "'40.65"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.154"
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0485"
$ws.Range("E42").Value = "  +0.22%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'3.13"
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.79"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.53"
$ws.Range("E45").Value = "  +4.15%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.145"
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'3.19"
$ws.Range("E47").Value = "  +2.41%  "
$ws.Range("B48").Value = "FLOKI"
$ws.Range("C48").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D48").Value = "'0.000282"
$ws.Range("E48").Value = "  +25.02%  "
$ws.Range("D49").Value = "'9.15"
$ws.Range("E49").Value = "  +6.88%  "
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").Value = "'3.33"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0344"
$ws.Range("E51").Value = "  +0.90%  "
